# Adds row 73 (match index 72) to the Armenia Premier League 2023-2024 sheet,
# mirroring the formatting of the last existing data row (row 72).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 73
$lastRow = 72

# --- Copy formatting from the previous data row, cell by cell, so that the
#     new row reuses the exact same style indices (bold/border style for
#     column A, date-time number format for column E, default for the rest).
$srcRange = $ws.Range("A" + $lastRow + ":V" + $lastRow)
$dstRange = $ws.Range("A" + $newRow + ":V" + $newRow)
$srcRange.Copy()
$dstRange.PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Fill in the values for the new row ---
$ws.Cells.Item($newRow, 1).Value = 72
$ws.Cells.Item($newRow, 2).Value = "armenia"
$ws.Cells.Item($newRow, 3).Value = "premier-league"
$ws.Cells.Item($newRow, 4).Value = "2023-2024"
$ws.Cells.Item($newRow, 5).Value = 45236.66666666666
$ws.Cells.Item($newRow, 6).Value = "BKMA"
$ws.Cells.Item($newRow, 7).Value = 1
$ws.Cells.Item($newRow, 8).Value = "Pyunik Yerevan"
$ws.Cells.Item($newRow, 9).Value = 4
$ws.Cells.Item($newRow, 10).Value = 9.949999999999999
$ws.Cells.Item($newRow, 11).Value = "05/11/2023 04:12"
$ws.Cells.Item($newRow, 12).Value = 11.32
$ws.Cells.Item($newRow, 13).Value = "06/11/2023 15:59"
$ws.Cells.Item($newRow, 14).Value = 6.58
$ws.Cells.Item($newRow, 15).Value = "05/11/2023 04:12"
$ws.Cells.Item($newRow, 16).Value = 6.71
$ws.Cells.Item($newRow, 17).Value = "06/11/2023 15:59"
$ws.Cells.Item($newRow, 18).Value = 1.18
$ws.Cells.Item($newRow, 19).Value = "05/11/2023 04:12"
$ws.Cells.Item($newRow, 20).Value = 1.22
$ws.Cells.Item($newRow, 21).Value = "06/11/2023 15:30"
$ws.Cells.Item($newRow, 22).Value = "https://www.betexplorer.com/football/armenia/premier-league/bkma-pyunik-yerevan/SlicuDZb/"
